# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.996.47"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.123.91"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'580.20"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'172.54"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "'37.16"
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").Value = "3.640.74"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "66.962.36"
$ws.Range("D17").Value = "'7.15"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "3.123.34"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "'476.08"
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").Value = "'0.709"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "'7.86"
$ws.Range("E22").Value = "  +5.13%  "
$ws.Range("D23").Value = "'83.88"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'13.23"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("D26").Value = "'10.32"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'7.93"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "'28.61"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  -6.37%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("D37").Value = "'46.96"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").Value = "'50.18"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").Value = "'382.39"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").Value = "'2.56"
$ws.Range("E46").Value = "  -9.15%  "
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D49").Value = "'24.98"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("E51").Value = "  -0.64%  "
